$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.584.11"
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = "  -0.24%  "

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.286.49"
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = "  -0.77%  "

$ws.Range("E4").Value = "  -0.03%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "95.94"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +3.55%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.28"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -1.22%  "

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  -0.93%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("E9").Value = "  -1.71%  "

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.51"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +1.76%  "

$ws.Range("E11").Value = "  +0.26%  "

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.99"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  -1.54%  "

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +0.73%  "

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.629.57"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -0.79%  "

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.28"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  +0.25%  "

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.847"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  -0.51%  "

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.287.30"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -1.55%  "

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.561.73"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -0.37%  "

$ws.Range("E19").Value = "  +2.60%  "

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.20"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -0.82%  "

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.15"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +0.88%  "

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.59"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +13.05%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.63"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -2.74%  "

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.08"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -6.21%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  +1.69%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.21"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -1.31%  "

$ws.Range("E28").Value = "  +2.53%  "

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.97"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +2.39%  "

$ws.Range("E30").Value = "  -6.28%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.04"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +1.98%  "

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.81"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -3.62%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0898"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -0.23%  "

$ws.Range("E34").Value = "  -4.07%  "

$ws.Range("E35").Value = "  -0.58%  "

$ws.Range("E36").Value = "  -2.26%  "

$ws.Range("E37").Value = "  -0.68%  "

$ws.Range("E38").Value = "  -2.43%  "

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.32"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -3.22%  "

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.242"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +3.53%  "

$ws.Range("E41").Value = "  +0.07%  "

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.30"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +1.21%  "

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.75"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +6.39%  "

$ws.Range("E44").Value = "  +1.62%  "

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.77"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -1.58%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.17"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -4.93%  "

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.102"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -0.71%  "

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "96.91"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -3.31%  "

$ws.Range("E49").Value = "  -1.05%  "

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.186"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +8.18%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.510.87"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -0.69%  "
